# Apply the updated cryptocurrency market data (price + 1h volume change,
# and in two cases a coin re-ranking) to Sheet1, cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.465.48"
$ws.Range("E2").Value = "  -0.33%  "
# Row 3
$ws.Range("D3").Value = "3.773.07"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.38"
$ws.Range("E5").Value = "  +0.27%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.17"
$ws.Range("E6").Value = "  -1.62%  "
# Row 7
$ws.Range("E7").Value = "  +0.04%  "
# Row 8
$ws.Range("E8").Value = "  -1.07%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.158"
$ws.Range("E9").Value = "  -1.15%  "
# Row 10
$ws.Range("E10").Value = "  +0.29%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").Value = "  +0.92%  "
# Row 12
$ws.Range("E12").Value = "  -2.17%  "
# Row 13
$ws.Range("E13").Value = "  -1.72%  "
# Row 14
$ws.Range("D14").Value = "4.411.53"
$ws.Range("E14").Value = "  -0.35%  "
# Row 15
$ws.Range("D15").Value = "3.759.06"
$ws.Range("E15").Value = "  -0.80%  "
# Row 16
$ws.Range("D16").Value = "67.537.95"
$ws.Range("E16").Value = "  -0.15%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.22"
$ws.Range("E17").Value = "  -1.70%  "
# Row 18
$ws.Range("E18").Value = "  +1.73%  "
# Row 19
$ws.Range("E19").Value = "  -0.27%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.17"
$ws.Range("E20").Value = "  -0.02%  "
# Row 21
$ws.Range("E21").Value = "  -3.07%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.692"
$ws.Range("E22").Value = "  -0.74%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000146"
$ws.Range("E23").Value = "  -5.59%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.28"
$ws.Range("E24").Value = "  -1.42%  "
# Row 25
$ws.Range("E25").Value = "  -1.72%  "
# Row 26
$ws.Range("E26").Value = "  -1.23%  "
# Row 27
$ws.Range("E27").Value = "  -0.10%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("E28").Value = "  -0.59%  "
# Row 29
$ws.Range("D29").Value = "3.925.09"
$ws.Range("E29").Value = "  -0.36%  "
# Row 30
$ws.Range("E30").Value = "  +3.06%  "
# Row 31
$ws.Range("E31").Value = "  -6.14%  "
# Row 32
$ws.Range("E32").Value = "  -2.81%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.00"
$ws.Range("E33").Value = "  -2.15%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.03%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.94"
$ws.Range("E35").Value = "  -1.29%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0984"
$ws.Range("E36").Value = "  -1.45%  "
# Row 37
$ws.Range("E37").Value = "  +0.30%  "
# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.22"
$ws.Range("E38").Value = "  -4.34%  "
# Row 39
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.985"
$ws.Range("E39").Value = "  -0.87%  "
# Row 40
$ws.Range("E40").Value = "  -0.68%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.06%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "47.36"
$ws.Range("E43").Value = "  -1.42%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.47"
$ws.Range("E44").Value = "  -0.72%  "
# Row 45
$ws.Range("E45").Value = "  -0.56%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.44"
$ws.Range("E46").Value = "  +1.04%  "
# Row 47
$ws.Range("E47").Value = "  +0.33%  "
# Row 48
$ws.Range("E48").Value = "  +8.48%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "27.13"
$ws.Range("E49").Value = "  +0.34%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.84"
$ws.Range("E50").Value = "  +1.07%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "390.10"
$ws.Range("E51").Value = "  +0.03%  "
